# NYPD CompStat weekly crime report refresh: new crime data collected.
# Updates the report week header text and the crime statistics table
# (104th Precinct) for the week of 1/22/2024 - 1/28/2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich text cells; runs share identical formatting
# so a plain value replacement renders identically).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# ---------------------------------------------------------------------
# Row 16 - Burglary
# ---------------------------------------------------------------------
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -57.142857142857
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -29.411764705882
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 17
$ws.Range("K16").Value = -29.411764705882
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = -42.857142857142
$ws.Range("N16").Value = -88.235294117647

# ---------------------------------------------------------------------
# Row 17 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 15.789473684210
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 19
$ws.Range("K17").Value = 15.789473684210
$ws.Range("L17").Value = 22.222222222222
$ws.Range("M17").Value = 29.411764705882
$ws.Range("N17").Value = 37.5

# ---------------------------------------------------------------------
# Row 18 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -38.888888888888
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = -38.888888888888
$ws.Range("L18").Value = -42.105263157894
$ws.Range("M18").Value = -73.170731707317
$ws.Range("N18").Value = -93.452380952380

# ---------------------------------------------------------------------
# Row 19 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 91.176470588235
$ws.Range("I19").Value = 65
$ws.Range("J19").Value = 34
$ws.Range("K19").Value = 91.176470588235
$ws.Range("L19").Value = 32.653061224489
$ws.Range("M19").Value = 75.675675675675
$ws.Range("N19").Value = 54.761904761904

# ---------------------------------------------------------------------
# Row 20 - Transit
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 21
$ws.Range("L20").Value = 21.739130434782
$ws.Range("M20").Value = 27.272727272727
$ws.Range("N20").Value = -90.378006872852

# ---------------------------------------------------------------------
# Row 21 - Housing (bold)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 31.034482758620
$ws.Range("F21").Value = 138
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = 25.454545454545
$ws.Range("I21").Value = 138
$ws.Range("J21").Value = 110
$ws.Range("K21").Value = 25.454545454545
$ws.Range("L21").Value = 15
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -77.813504823151

# ---------------------------------------------------------------------
# Row 22 - Petit Larceny
# D22/E22 switch from the " "/"***.*" placeholder text to real numbers,
# so the number format has to be (re)applied explicitly to land on the
# same styles used by the other numeric cells in the table.
# ---------------------------------------------------------------------
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 0

# ---------------------------------------------------------------------
# Row 24 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -25.806451612903
$ws.Range("F24").Value = 116
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = 2.654867256637
$ws.Range("I24").Value = 116
$ws.Range("J24").Value = 113
$ws.Range("K24").Value = 2.654867256637
$ws.Range("L24").Value = 18.367346938775
$ws.Range("M24").Value = 34.883720930232

# ---------------------------------------------------------------------
# Row 25 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -27.272727272727
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -17.142857142857
$ws.Range("I25").Value = 29
$ws.Range("J25").Value = 35
$ws.Range("K25").Value = -17.142857142857
$ws.Range("L25").Value = -6.451612903225
$ws.Range("M25").Value = -39.583333333333

# ---------------------------------------------------------------------
# Row 27 - Hate Crimes
# C27 switches from a real number back to the " " placeholder text.
# Number format is toggled through "@" (text) so the numeric-looking
# string "0" is stored as text instead of being re-parsed as a number,
# then restored to general so it lands back on the shared text style.
# ---------------------------------------------------------------------
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C27").NumberFormat = "general"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -66.666666666666
$ws.Range("L27").Value = -50
